# "案件情報.xlsx" append: 2025-11-09 18:23 JST
# Three freshly-scraped Lancers listings are prepended to the "ランサーズ" sheet
# (rows 2-4), pushing the previously-captured rows down by three, and the capture
# timestamp in column A is refreshed to 2025-11-09 18:23:48 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room at the top: insert 3 blank rows before the current row 2, which
# shifts all existing data rows (old 2..12) down to (new 5..15).
$ws.Rows("2:4").Insert()

# The row insert does not carry the Hyperlinks collection along with it (the
# anchors stay pinned to F2:F12), so drop the stale links; correct ones for
# every data row (new layout + the 3 new rows) are (re)created below.
$ws.Hyperlinks.Delete()

# Row 2  (new)
$ws.Cells.Item(2, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(2, 2).Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5423720")
$ws.Cells.Item(2, 7).Value = 385
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆効率化"

# Row 3  (new)
$ws.Cells.Item(3, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(3, 2).Value = "【募集】習慣化+目標管理を目的としたAIネイティブなWebサービスのMVP開発"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5430365"
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5430365")
$ws.Cells.Item(3, 7).Value = 378
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◆開発 ◇管理"

# Row 4  (new)
$ws.Cells.Item(4, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(4, 2).Value = "製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5419380"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5419380")
$ws.Cells.Item(4, 7).Value = 298
$ws.Cells.Item(4, 8).Value = "🔥AI,Ai"

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(5, 2).Value = "海外仕入れ元サイト→ツールを動かす為のCSVファイルに週1で自動抽出の制作(自動/スクレイピング)"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5251319"
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5251319")
$ws.Cells.Item(5, 7).Value = 135
$ws.Cells.Item(5, 8).Value = "◆ツール,スクレイピング ◇サイト"

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(6, 2).Value = "【急募】Wordpressを用いた比較サイトの新規開発"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5430121"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5430121")
$ws.Cells.Item(6, 7).Value = 123
$ws.Cells.Item(6, 8).Value = "◆開発 ◇サイト ○WordPress"

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(7, 2).Value = "Glideメインで作成したシステムをLinux+MySQL型に移行するための新規開発"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5430095"
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5430095")
$ws.Cells.Item(7, 7).Value = 115
$ws.Cells.Item(7, 8).Value = "◆開発 ◇MySQL"

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(8, 2).Value = "初回 webアプリの開発"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5430337"
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5430337")
$ws.Cells.Item(8, 7).Value = 100
$ws.Cells.Item(8, 8).Value = "◆開発 ◇アプリ"

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(9, 2).Value = "【急募】ビデオサーバー開発:MXFファイル再生とHD-SDI出力"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5430205"
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://www.lancers.jp/work/detail/5430205")
$ws.Cells.Item(9, 7).Value = 75
$ws.Cells.Item(9, 8).Value = "◆開発"

# Row 10
$ws.Cells.Item(10, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(10, 2).Value = "小売店向けシステム性能試験"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5430176"
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://www.lancers.jp/work/detail/5430176")
$ws.Cells.Item(10, 7).Value = 40
$ws.Cells.Item(10, 8).ClearContents()

# Row 11
$ws.Cells.Item(11, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(11, 2).Value = "【急募】YAMAHA RTX1210のVLAN設定作業"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5430327"
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://www.lancers.jp/work/detail/5430327")
$ws.Cells.Item(11, 7).Value = 13
$ws.Cells.Item(11, 8).ClearContents()

# Row 12
$ws.Cells.Item(12, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(12, 2).Value = "初回 飲食手向けハンディ、モバイルオーダー検証"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5430301"
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), "https://www.lancers.jp/work/detail/5430301")
$ws.Cells.Item(12, 7).Value = 13
$ws.Cells.Item(12, 8).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(13, 2).Value = "インターネット情報収集(selenium)"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5430171"
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://www.lancers.jp/work/detail/5430171")
$ws.Cells.Item(13, 7).Value = 10
$ws.Cells.Item(13, 8).ClearContents()

# Row 14
$ws.Cells.Item(14, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(14, 2).Value = "MT4 RSXを使用したEAの作成依頼"
$ws.Cells.Item(14, 3).Value = "システム開発"
$ws.Cells.Item(14, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(14, 5).Value = "期限情報なし"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5430008"
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), "https://www.lancers.jp/work/detail/5430008")
$ws.Cells.Item(14, 7).Value = 10
$ws.Cells.Item(14, 8).ClearContents()

# Row 15
$ws.Cells.Item(15, 1).Value = "2025-11-09 18:23:48"
$ws.Cells.Item(15, 2).Value = "【急募】LINE × QRコード連携で自動取得設定を実現!"
$ws.Cells.Item(15, 3).Value = "システム開発"
$ws.Cells.Item(15, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(15, 5).Value = "期限情報なし"
$ws.Cells.Item(15, 6).Value = "https://www.lancers.jp/work/detail/5430015"
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), "https://www.lancers.jp/work/detail/5430015")
$ws.Cells.Item(15, 7).Value = 10
$ws.Cells.Item(15, 8).ClearContents()
